# Applies the OrangeHRM_Excel.xlsx commit:
#  - Admin sheet: D2 value changes to "7598587A"
#  - Jira sheet: B2 cleared, A3/B3 updated, and three new rows (4-6)
#    appended with Recruitment_HiredList / Recruitment_RejectionList1 /
#    PersonalDetails in column A (column B left blank on every new row).

$wb = $excel.ActiveWorkbook

# --- Admin sheet --------------------------------------------------------
$admin = $wb.Worksheets.Item("Admin")
$admin.Range("D2").Value = "7598587A"

# --- Jira sheet ----------------------------------------------------------
$jira = $wb.Worksheets.Item("Jira")

# Row 2: Key / (value cleared)
$jira.Range("B2").Value = ""

# Row 3: id column now holds "PersonalDetails", value column cleared (but
# still present as a (blank) cell)
$jira.Range("A3").Value = "PersonalDetails"
$jira.Range("B3").Value = ""
$jira.Range("B3").Style = "Normal"

# New rows 4-6 - column A gets the new id, column B stays blank but present
$jira.Range("A4").Value = "Recruitment_HiredList"
$jira.Range("B4").Style = "Normal"

$jira.Range("A5").Value = "Recruitment_RejectionList1"
$jira.Range("B5").Style = "Normal"

$jira.Range("A6").Value = "PersonalDetails"
$jira.Range("B6").Style = "Normal"
